# Memory map: add "ELF file loader" entry (first pass of ELF definitions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at 14 (pushes old rows 14-16 down to 15-17),
#    shifting formulas automatically.
# ------------------------------------------------------------------
$ws.Rows.Item(14).Insert()

# ------------------------------------------------------------------
# 2. Row 13 ("Shell Code and Static Data") now ends where the new
#    ELF loader entry begins, so its size changes from E00 to 1C00.
# ------------------------------------------------------------------
$ws.Range("E13").Formula = '=T("1C00")'

# ------------------------------------------------------------------
# 3. Fill in the new row 14 - "ELF file loader" (ELFL.RNB)
# ------------------------------------------------------------------
$ws.Range("A14").Formula = "=F13"
$ws.Range("B14").Value = "'ELF file loader"
$ws.Range("C14").Value = "'ELFL.RNB"
$ws.Range("D14").Value = "'True"
$ws.Range("E14").Formula = '=T("1C00")'
$ws.Range("F14").Formula = "=DEC2HEX((HEX2DEC(A14)+HEX2DEC(E14)),10)"

# ------------------------------------------------------------------
# 4. Give F14 its own distinct (18th) font - plain Arial 10, explicit
#    black - so it carries a dedicated style entry like in the source.
# ------------------------------------------------------------------
$ws.Range("F14").Font.Color = 0

# ------------------------------------------------------------------
# 5. Recolor the "Free Memory" rows (9, 12 and the row that used to
#    be 14, now 15) with the new blue fill (fgColor 007FFF / bgColor
#    3366FF), matching the palette color used for "EBDA"/"Video
#    Memory" entries.
# ------------------------------------------------------------------
foreach ($rowNum in 9, 12, 15) {
    $rowRange = $ws.Range("A" + $rowNum + ":F" + $rowNum)
    $rowRange.Interior.Color = 16744192
    $rowRange.Interior.PatternColor = 16737843
}

# ------------------------------------------------------------------
# 6. Restore the scientific-notation number format on column A/E of
#    those rows where it belongs (A9/A12 stay General; only the
#    original "A14"-now-"A15" cell keeps the 0.00E+00 custom format,
#    which the row-insert already preserved automatically).
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# 7. Move the active selection to E15 (mirrors the author's cursor
#    position after editing row 13's size formula that now lives one
#    row lower in the sheet).
# ------------------------------------------------------------------
$ws.Range("E15").Select()

Write-Host "edit applied"
